$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns.
# A leading apostrophe forces text entry so values stay strings
# (matching the source cells original inline-string type) instead of
# being auto-converted to numbers by Excel.
$ws.Range('D2').Value = "'40.036.77"
$ws.Range('E2').Value = "'  +0.58%  "
$ws.Range('D3').Value = "'2.211.90"
$ws.Range('E3').Value = "'  -0.55%  "
$ws.Range('E4').Value = "'  -0.05%  "
$ws.Range('D5').Value = "'290.24"
$ws.Range('E5').Value = "'  -3.12%  "
$ws.Range('D6').Value = "'88.51"
$ws.Range('E6').Value = "'  +4.83%  "
$ws.Range('E7').Value = "'  +0.36%  "
$ws.Range('E8').Value = "'  -0.10%  "
$ws.Range('D9').Value = "'0.471"
$ws.Range('E9').Value = "'  +1.00%  "
$ws.Range('D10').Value = "'30.78"
$ws.Range('E10').Value = "'  +3.44%  "
$ws.Range('E11').Value = "'  +0.20%  "
$ws.Range('D12').Value = "'47.65"
$ws.Range('E12').Value = "'  +2.14%  "
$ws.Range('E13').Value = "'  +2.13%  "
$ws.Range('D14').Value = "'6.45"
$ws.Range('E14').Value = "'  +2.43%  "
$ws.Range('D15').Value = "'2.556.20"
$ws.Range('E15').Value = "'  -0.76%  "
$ws.Range('D16').Value = "'14.00"
$ws.Range('E16').Value = "'  -0.77%  "
$ws.Range('D17').Value = "'2.221.83"
$ws.Range('E17').Value = "'  -0.13%  "
$ws.Range('E18').Value = "'  +1.32%  "
$ws.Range('D19').Value = "'39.964.34"
$ws.Range('E19').Value = "'  +0.52%  "
$ws.Range('D20').Value = "'11.78"
$ws.Range('E20').Value = "'  +13.02%  "
$ws.Range('E21').Value = "'  +0.95%  "
$ws.Range('D22').Value = "'5.80"
$ws.Range('E22').Value = "'  +0.73%  "
$ws.Range('D23').Value = "'65.73"
$ws.Range('E23').Value = "'  +1.06%  "
$ws.Range('D24').Value = "'236.20"
$ws.Range('E24').Value = "'  +0.77%  "
$ws.Range('E25').Value = "'  +0.05%  "
$ws.Range('D26').Value = "'2.47"
$ws.Range('E26').Value = "'  +1.47%  "
$ws.Range('E27').Value = "'  +1.45%  "
$ws.Range('D28').Value = "'22.61"
$ws.Range('E28').Value = "'  -0.63%  "
$ws.Range('D29').Value = "'2.19"
$ws.Range('E29').Value = "'  +4.10%  "
$ws.Range('E30').Value = "'  +0.60%  "
$ws.Range('D31').Value = "'153.31"
$ws.Range('E31').Value = "'  +2.26%  "
$ws.Range('D32').Value = "'32.16"
$ws.Range('E32').Value = "'  -0.65%  "
$ws.Range('E33').Value = "'  -0.12%  "
$ws.Range('E34').Value = "'  +2.67%  "
$ws.Range('D35').Value = "'0.0718"
$ws.Range('E35').Value = "'  +2.39%  "
$ws.Range('E36').Value = "'  -0.51%  "
$ws.Range('E37').Value = "'  +7.05%  "
$ws.Range('D38').Value = "'16.04"
$ws.Range('E38').Value = "'  -2.24%  "
$ws.Range('D39').Value = "'0.111"
$ws.Range('E39').Value = "'  +0.68%  "
$ws.Range('D40').Value = "'0.100"
$ws.Range('E40').Value = "'  +2.56%  "
$ws.Range('E41').Value = "'  +2.69%  "
$ws.Range('D42').Value = "'2.097.47"
$ws.Range('E42').Value = "'  +8.60%  "
$ws.Range('E43').Value = "'  +4.42%  "
$ws.Range('E44').Value = "'  +2.30%  "
$ws.Range('D45').Value = "'0.0269"
$ws.Range('E45').Value = "'  +1.39%  "
$ws.Range('D46').Value = "'9.88"
$ws.Range('E46').Value = "'  +6.98%  "
$ws.Range('D47').Value = "'17.67"
$ws.Range('E47').Value = "'  +6.94%  "
$ws.Range('E48').Value = "'  +1.98%  "
$ws.Range('D49').Value = "'2.429.93"
$ws.Range('E49').Value = "'  -0.64%  "
$ws.Range('D50').Value = "'69.71"
$ws.Range('E50').Value = "'  -1.77%  "
$ws.Range('D51').Value = "'88.64"
$ws.Range('E51').Value = "'  +0.09%  "
